# Auto-generated Excel COM-interop script to apply F-column ('想去人数') updates
# across sheets 展览, 演出, and 全部类型, per the provided diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1030
$ws.Range("F4").Value = 5805
$ws.Range("F5").Value = 541
$ws.Range("F6").Value = 1011
$ws.Range("F7").Value = 1023
$ws.Range("F8").Value = 842
$ws.Range("F10").Value = 47
$ws.Range("F11").Value = 611
$ws.Range("F15").Value = 1972
$ws.Range("F16").Value = 1506
$ws.Range("F17").Value = 1056
$ws.Range("F19").Value = 206
$ws.Range("F20").Value = 385
$ws.Range("F21").Value = 620
$ws.Range("F22").Value = 220
$ws.Range("F23").Value = 1067
$ws.Range("F26").Value = 3417
$ws.Range("F27").Value = 191
$ws.Range("F28").Value = 127
$ws.Range("F30").Value = 151
$ws.Range("F31").Value = 49
$ws.Range("F32").Value = 478
$ws.Range("F33").Value = 15
$ws.Range("F36").Value = 227
$ws.Range("F37").Value = 316
$ws.Range("F38").Value = 805
$ws.Range("F40").Value = 68
$ws.Range("F41").Value = 74
$ws.Range("F42").Value = 82

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 532
$ws.Range("F6").Value = 294

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1030
$ws.Range("F6").Value = 5805
$ws.Range("F7").Value = 541
$ws.Range("F8").Value = 1011
$ws.Range("F10").Value = 532
$ws.Range("F11").Value = 1023
$ws.Range("F12").Value = 842
$ws.Range("F14").Value = 294
$ws.Range("F16").Value = 47
$ws.Range("F17").Value = 611
$ws.Range("F22").Value = 1972
$ws.Range("F23").Value = 1506
$ws.Range("F24").Value = 1056
$ws.Range("F26").Value = 206
$ws.Range("F27").Value = 385
$ws.Range("F29").Value = 620
$ws.Range("F30").Value = 220
$ws.Range("F31").Value = 1067
$ws.Range("F32").Value = 3417
$ws.Range("F33").Value = 191
$ws.Range("F34").Value = 127
$ws.Range("F36").Value = 151
$ws.Range("F37").Value = 49
$ws.Range("F38").Value = 478
$ws.Range("F39").Value = 15
$ws.Range("F42").Value = 316
$ws.Range("F43").Value = 805
$ws.Range("F45").Value = 74
$ws.Range("F46").Value = 82
